# Scheduled-runner price/profit refresh across job-leve sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) for rows whose underlying market data changed; some rows
# also gain/lose the NQ-only (M) vs HQ (N) profit cell depending on
# whether an HQ price is available.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""

$ws.Range("H138").Value = 1820785.5
$ws.Range("I138").Value = 2633441
$ws.Range("J138").Value = 4261.4707
$ws.Range("K138").Value = 7900323
$ws.Range("L138").Value = 12784.4121
$ws.Range("M138").Value = -7895183
$ws.Range("N138").Value = -23064.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5298.3037
$ws.Range("I32").Value = 5722.7954
$ws.Range("J32").Value = 3741.8333
$ws.Range("K32").Value = 5722.7954
$ws.Range("L32").Value = 3741.8333
$ws.Range("M32").Value = -5435.7954
$ws.Range("N32").Value = -4315.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 49890
$ws.Range("J132").Value = 49890
$ws.Range("L132").Value = 49890
$ws.Range("N132").Value = -60010

$ws.Range("H134").Value = 23817.137
$ws.Range("I134").Value = 26663.719
$ws.Range("K134").Value = 79991.15700000001
$ws.Range("M134").Value = -77456.15700000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 963.5599999999999
$ws.Range("I5").Value = 848.1667
$ws.Range("K5").Value = 2544.5001
$ws.Range("M5").Value = -2432.5001

$ws.Range("H63").Value = 2250
$ws.Range("I63").Value = 2250
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 6750
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -6001
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 2250
$ws.Range("I66").Value = 2250
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 20250
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -16506
$ws.Range("N66").Value = ""

$ws.Range("H76").Value = 1325
$ws.Range("J76").Value = 1962.5
$ws.Range("L76").Value = 5887.5
$ws.Range("N76").Value = -6653.5

$ws.Range("H79").Value = 1325
$ws.Range("J79").Value = 1962.5
$ws.Range("L79").Value = 5887.5
$ws.Range("N79").Value = -8539.5

$ws.Range("H81").Value = 1693.75
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1693.75
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 5081.25
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = -7327.25

$ws.Range("H84").Value = 1693.75
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1693.75
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 15243.75
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = -26475.75

$ws.Range("H94").Value = 3541
$ws.Range("J94").Value = 4095.25
$ws.Range("L94").Value = 12285.75
$ws.Range("N94").Value = -13637.75

$ws.Range("H96").Value = 35356104
$ws.Range("J96").Value = 35356104
$ws.Range("L96").Value = 106068312
$ws.Range("N96").Value = -106072430

$ws.Range("H101").Value = 6225.625
$ws.Range("J101").Value = 6939.857
$ws.Range("L101").Value = 20819.571
$ws.Range("N101").Value = -25687.571

$ws.Range("H106").Value = 3972998.2
$ws.Range("J106").Value = 3972998.2
$ws.Range("L106").Value = 11918994.6
$ws.Range("N106").Value = -11920886.6

$ws.Range("H108").Value = 1149.4
$ws.Range("I108").Value = 1149.4
$ws.Range("K108").Value = 3448.2
$ws.Range("M108").Value = -568.2000000000003

$ws.Range("H110").Value = 3600
$ws.Range("J110").Value = 3620
$ws.Range("L110").Value = 10860
$ws.Range("N110").Value = -19040

$ws.Range("H113").Value = 699.1539
$ws.Range("I113").Value = 730
$ws.Range("J113").Value = 689.9
$ws.Range("K113").Value = 2190
$ws.Range("L113").Value = 2069.7
$ws.Range("M113").Value = -20
$ws.Range("N113").Value = -6409.7

$ws.Range("H122").Value = 22224346
$ws.Range("I122").Value = 47619776
$ws.Range("J122").Value = 3343.75
$ws.Range("K122").Value = 428577984
$ws.Range("L122").Value = 30093.75
$ws.Range("M122").Value = -428575534
$ws.Range("N122").Value = -34993.75

$ws.Range("H135").Value = 963.5599999999999
$ws.Range("I135").Value = 848.1667
$ws.Range("K135").Value = 7633.5003
$ws.Range("M135").Value = -5098.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5331.25
$ws.Range("I70").Value = 5330
$ws.Range("J70").Value = 5333.3335
$ws.Range("K70").Value = 5330
$ws.Range("L70").Value = 5333.3335
$ws.Range("M70").Value = -5060
$ws.Range("N70").Value = -5873.3335

$ws.Range("H73").Value = 5331.25
$ws.Range("I73").Value = 5330
$ws.Range("J73").Value = 5333.3335
$ws.Range("K73").Value = 5330
$ws.Range("L73").Value = 5333.3335
$ws.Range("M73").Value = -4394
$ws.Range("N73").Value = -7205.3335

$ws.Range("H80").Value = 9362.058999999999
$ws.Range("I80").Value = 3259.1667
$ws.Range("J80").Value = 12690.909
$ws.Range("K80").Value = 3259.1667
$ws.Range("L80").Value = 12690.909
$ws.Range("M80").Value = -2261.1667
$ws.Range("N80").Value = -14686.909

$ws.Range("H83").Value = 9362.058999999999
$ws.Range("I83").Value = 3259.1667
$ws.Range("J83").Value = 12690.909
$ws.Range("K83").Value = 16295.8335
$ws.Range("L83").Value = 63454.545
$ws.Range("M83").Value = -11303.8335
$ws.Range("N83").Value = -73438.545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2677.7144
$ws.Range("I7").Value = 2584
$ws.Range("J7").Value = 2771.4285
$ws.Range("K7").Value = 2584
$ws.Range("L7").Value = 2771.4285
$ws.Range("M7").Value = -2472
$ws.Range("N7").Value = -2995.4285

$ws.Range("H126").Value = 2677.7144
$ws.Range("I126").Value = 2584
$ws.Range("J126").Value = 2771.4285
$ws.Range("K126").Value = 7752
$ws.Range("L126").Value = 8314.2855
$ws.Range("M126").Value = -5282
$ws.Range("N126").Value = -13254.2855

$ws.Range("H132").Value = 3513.0637
$ws.Range("I132").Value = 3509.5945
$ws.Range("J132").Value = 3525.9
$ws.Range("K132").Value = 10528.7835
$ws.Range("L132").Value = 10577.7
$ws.Range("M132").Value = -7998.783500000001
$ws.Range("N132").Value = -15637.7

$ws.Range("H136").Value = 2214.0625
$ws.Range("I136").Value = 1262.9166
$ws.Range("K136").Value = 3788.7498
$ws.Range("M136").Value = -1238.7498
